$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clothing = @{
    2 = 'Jumpsuit,Dress'
    3 = 'Halter,Parka'
    4 = 'Jumpsuit,Kaftan'
    5 = 'Parka,Jumpsuit'
    6 = 'Halter,Parka'
    7 = 'Parka,Jumpsuit'
    8 = 'Jumpsuit,Kaftan'
    9 = 'Halter,Tee'
    10 = 'Blouse,Jumpsuit'
    11 = 'Jumpsuit,Blouse'
    12 = 'Parka,Cutoffs'
    13 = 'Trunks,Coverup'
    14 = 'Blouse,Jumpsuit'
    15 = 'Jumpsuit,Kaftan'
    16 = 'Blazer,Trunks'
    17 = 'Caftan,Trunks'
    18 = 'Blouse,Caftan'
    19 = 'Jumpsuit,Kaftan'
    20 = 'Parka,Gauchos'
    21 = 'Jumpsuit,Kaftan'
    22 = 'Jumpsuit,Blouse'
    23 = 'Jumpsuit,Dress'
    24 = 'Jumpsuit,Blazer'
    25 = 'Blouse,Trunks'
    26 = 'Jumpsuit,Kaftan'
    27 = 'Kaftan,Jumpsuit'
    28 = 'Jumpsuit,Dress'
    29 = 'Dress,Blouse'
    30 = 'Jumpsuit,Tee'
    31 = 'Halter,Parka'
    32 = 'Caftan,Jumpsuit'
    33 = 'Trunks,Caftan'
    34 = 'Blouse,Jumpsuit'
    35 = 'Jumpsuit,Dress'
    36 = 'Blouse,Jumpsuit'
    37 = 'Jumpsuit,Parka'
    38 = 'Jumpsuit,Halter'
    39 = 'Jumpsuit,Blazer'
    40 = 'Parka,Caftan'
    41 = 'Tee,Halter'
    42 = 'Blouse,Jumpsuit'
    43 = 'Tee,Blazer'
    44 = 'Sweatpants,Dress'
    45 = 'Jumpsuit,Kaftan'
    46 = 'Jumpsuit,Blouse'
    47 = 'Jodhpurs,Jumpsuit'
    48 = 'Trunks,Blouse'
    49 = 'Trunks,Jumpsuit'
    50 = 'Trunks,Jodhpurs'
    51 = 'Parka,Caftan'
    52 = 'Jumpsuit,Cutoffs'
    53 = 'Blouse,Blazer'
    54 = 'Blouse,Jumpsuit'
    55 = 'Jumpsuit,Blazer'
    56 = 'Jumpsuit,Dress'
    57 = 'Jumpsuit,Blouse'
    58 = 'Jumpsuit,Caftan'
    59 = 'Jumpsuit,Tee'
    60 = 'Jumpsuit,Halter'
    61 = 'Halter,Blazer'
    62 = 'Jumpsuit,Chinos'
    63 = 'Jumpsuit,Kaftan'
    64 = 'Jumpsuit,Blouse'
    65 = 'Blouse,Parka'
    66 = 'Trunks,Halter'
    67 = 'Blouse,Jumpsuit'
    68 = 'Parka,Jumpsuit'
    69 = 'Jumpsuit,Kaftan'
    70 = 'Jumpsuit,Blouse'
    71 = 'Blouse,Caftan'
    72 = 'Kaftan,Jumpsuit'
    73 = 'Jumpsuit,Dress'
    74 = 'Blouse,Jumpsuit'
    75 = 'Jumpsuit,Kaftan'
    76 = 'Parka,Caftan'
    77 = 'Trunks,Caftan'
    78 = 'Jumpsuit,Tee'
    79 = 'Parka,Jumpsuit'
    80 = 'Halter,Jumpsuit'
    81 = 'Jumpsuit,Blouse'
    82 = 'Blazer,Halter'
    83 = 'Halter,Blazer'
    84 = 'Blazer,Halter'
    85 = 'Halter,Jumpsuit'
    86 = 'Halter,Blouse'
    87 = 'Parka,Gauchos'
    88 = 'Parka,Trunks'
    89 = 'Halter,Jumpsuit'
    90 = 'Dress,Jumpsuit'
    91 = 'Jumpsuit,Blazer'
    92 = 'Blouse,Parka'
    93 = 'Jumpsuit,Dress'
    94 = 'Jumpsuit,Trunks'
}

foreach ($row in $clothing.Keys) {
    $ws.Range("G$row").Value = $clothing[$row]
}
